$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6021503210067749
$ws.Range("B1").Value = 0.9641004800796509
$ws.Range("C1").Value = 4.289170265197754
$ws.Range("D1").Value = 1.976656794548035
$ws.Range("E1").Value = 1.605364322662354
